# Update the player roster table on Sheet1 (columns A:C, rows 2-19)
# to reflect the new roster/positions/teams per the commit's data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Amen Thompson",      "SG,SF",    "Houston Rockets"),
    @("Jabari Smith Jr.",   "PF,C",     "Houston Rockets"),
    @("Yves Missi",         "C",        "New Orleans Pelicans"),
    @("Alperen Sengün",     "C",        "Houston Rockets"),
    @("Coby White",         "PG,SG",    "Chicago Bulls"),
    @("Kelly Oubre Jr.",    "SG,SF",    "Philadelphia 76ers"),
    @("Caris LeVert",       "SG,SF",    "Cleveland Cavaliers"),
    @("Norman Powell",      "SG,SF",    "LA Clippers"),
    @("Jeremy Sochan",      "SF,PF",    "San Antonio Spurs"),
    @("Devin Vassell",      "SG,SF",    "San Antonio Spurs"),
    @("Desmond Bane",       "SG,SF",    "Memphis Grizzlies"),
    @("Walker Kessler",     "C",        "Utah Jazz"),
    @("LeBron James",       "SF,PF",    "Los Angeles Lakers"),
    @("Jalen Brunson",      "PG",       "New York Knicks"),
    @("Trae Young",         "PG",       "Atlanta Hawks"),
    @("Immanuel Quickley",  "PG,SG",    "Toronto Raptors"),
    @("Kawhi Leonard",      "SG,SF,PF", "LA Clippers"),
    @("Devin Booker",       "PG,SG",    "Phoenix Suns")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
